# Apply COVID-19 country statistics update (commit: "Update countries & provincias Spain")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp shown in cell A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 20:21"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2344053
$ws.Range("C4").Value = 13475
$ws.Range("D4").Value = 975049
$ws.Range("E4").Value = 1246877
$ws.Range("G4").Value = 147
$ws.Range("H4").Value = 122127

# Row 7: India
$ws.Range("B7").Value = 426473
$ws.Range("C7").Value = 14746
$ws.Range("D7").Value = 236927
$ws.Range("E7").Value = 175851
$ws.Range("G7").Value = 418
$ws.Range("H7").Value = 13695

# Row 11: Chile
$ws.Range("A11").Value = "Chile"
$ws.Range("B11").Value = 242355
$ws.Range("C11").Value = 5607
$ws.Range("D11").Value = 200569
$ws.Range("E11").Value = 37307
$ws.Range("G11").Value = 184
$ws.Range("H11").Value = 4479

# Row 12: Italia
$ws.Range("A12").Value = "Italia"
$ws.Range("B12").Value = 238499
$ws.Range("C12").Value = 224
$ws.Range("D12").Value = 182893
$ws.Range("E12").Value = 20972
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 34634

# Row 14: Alemania
$ws.Range("B14").Value = 191321
$ws.Range("C14").Value = 105
$ws.Range("E14").Value = 7460

# Row 15: Turquia
$ws.Range("B15").Value = 187685
$ws.Range("C15").Value = 1192
$ws.Range("D15").Value = 160240
$ws.Range("E15").Value = 22495
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = 4950

# Row 51: Israel
$ws.Range("B51").Value = 20741
$ws.Range("C51").Value = 108
$ws.Range("D51").Value = 15689
$ws.Range("E51").Value = 4746

# Row 58: Ghana
$ws.Range("B58").Value = 14007
$ws.Range("C58").Value = 290
$ws.Range("D58").Value = 10473
$ws.Range("E58").Value = 3449

# Row 68: Marruecos
$ws.Range("B68").Value = 9977
$ws.Range("C68").Value = 138
$ws.Range("D68").Value = 8284
$ws.Range("E68").Value = 1479
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 214

# Row 103: Maldivas
$ws.Range("B103").Value = 2203
$ws.Range("C103").Value = 16
$ws.Range("D103").Value = 1803
$ws.Range("E103").Value = 392

# Row 109: Sudan del Sur
$ws.Range("B109").Value = 1892
$ws.Range("C109").Value = 10
$ws.Range("D109").Value = 169
$ws.Range("E109").Value = 1689

# Row 121: Paraguay
$ws.Range("B121").Value = 1379
$ws.Range("C121").Value = 17
$ws.Range("D121").Value = 871
$ws.Range("E121").Value = 495

# Row 132: Cabo Verde
$ws.Range("A132").Value = "Cabo Verde"
$ws.Range("B132").Value = 890
$ws.Range("C132").Value = 27
$ws.Range("D132").Value = 377
$ws.Range("E132").Value = 505
$ws.Range("H132").Value = 8

# Row 133: Congo
$ws.Range("A133").Value = "Congo"
$ws.Range("B133").Value = 883
$ws.Range("D133").Value = 391
$ws.Range("E133").Value = 465
$ws.Range("H133").Value = 27

# Row 137: Estado de Palestina
$ws.Range("B137").Value = 833
$ws.Range("C137").Value = 49
$ws.Range("E137").Value = 391

# Row 153: Reunion
$ws.Range("B153").Value = 506
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 45

# Row 156: Montenegro
$ws.Range("B156").Value = 362
$ws.Range("C156").Value = 3
$ws.Range("E156").Value = 38

# Row 202: Fiyi
$ws.Range("A202").Value = "Fiyi"

# Row 203: Dominica
$ws.Range("A203").Value = "Dominica"

# Row 207: Groenlandia
$ws.Range("A207").Value = "Groenlandia"

# Row 208: Islas Malvinas
$ws.Range("A208").Value = "Islas Malvinas"

# Row 213: Papua Nueva Guinea
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214: Islas Virgenes Britanicas
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
